# Generate Report for handback
# Adds a new handback entry (47f6ed0b-d486-4127-a6ac-09c1a1d88e98) as row 4
# on the "Overview", "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$uuid = "47f6ed0b-d486-4127-a6ac-09c1a1d88e98"
$xlfHash = "e3dcd7592d74dd8ca9acefa64bc5f0ac4f434d82"

$mdName = "$uuid.md"
$statusInSync = "Handed back: in sync with en-US"
$includeReason = "Include"

$zhXlfName = "$uuid.$xlfHash.zh-cn.xlf"
$deXlfName = "$uuid.$xlfHash.de-de.xlf"

$zhHandoffDt = "2016-01-19 04:14:52"
$zhHandbackDt = "2016-01-19 04:15:33"
$deHandoffDt = "2016-01-19 04:15:02"
$deHandbackDt = "2016-01-19 04:15:50"

# Hyperlink target URLs (follow the same github blob url convention
# used for the existing rows in each worksheet).
$zhMdUrlA = "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000a/e2e/$mdName"
$zhMdUrlE = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000c/e2e/$mdName"
$zhXlfUrlC = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName"
$zhXlfUrlF = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName"

$deMdUrlA = $zhMdUrlA
$deMdUrlE = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000f/e2e/$mdName"
$deXlfUrlC = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName"
$deXlfUrlF = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/00000000000000000000000000000000000010/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName"

function Style-AsHyperlink($cell) {
    # Match the existing custom "HyperLink" cell style used by the workbook
    # (underlined, cornflower-blue Calibri 11) rather than Excel's default
    # theme hyperlink style that .Hyperlinks.Add() applies automatically.
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $zhMdUrlA, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
Style-AsHyperlink $wsOverview.Range("A4")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Range("C4").Value = $zhXlfName
$wsZh.Range("D4").Value = $zhHandoffDt
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $zhXlfName
$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $includeReason

$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $zhMdUrlA, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
Style-AsHyperlink $wsZh.Range("A4")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhXlfUrlC, [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null
Style-AsHyperlink $wsZh.Range("C4")
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), $zhMdUrlE, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
Style-AsHyperlink $wsZh.Range("E4")
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), $zhXlfUrlF, [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null
Style-AsHyperlink $wsZh.Range("F4")

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Range("C4").Value = $deXlfName
$wsDe.Range("D4").Value = $deHandoffDt
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $deXlfName
$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $includeReason

$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $deMdUrlA, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
Style-AsHyperlink $wsDe.Range("A4")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deXlfUrlC, [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null
Style-AsHyperlink $wsDe.Range("C4")
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), $deMdUrlE, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
Style-AsHyperlink $wsDe.Range("E4")
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), $deXlfUrlF, [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null
Style-AsHyperlink $wsDe.Range("F4")

Write-Host "Added handback row for $uuid to Overview, zh-cn and de-de sheets."
